$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-12 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-13 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("25+68=93", $true, $false, $false, $false, $false, $true, 1, $false, "60+8=68", 2) | Out-Null
$d.Content.Find.Execute("67-19=48", $true, $false, $false, $false, $false, $true, 1, $false, "35-30=5", 2) | Out-Null
$d.Content.Find.Execute("37+29=66", $true, $false, $false, $false, $false, $true, 1, $false, "50+10=60", 2) | Out-Null
$d.Content.Find.Execute("14+32=46", $true, $false, $false, $false, $false, $true, 1, $false, "19+19=38", 2) | Out-Null
$d.Content.Find.Execute("5+76=81", $true, $false, $false, $false, $false, $true, 1, $false, "71-37=34", 2) | Out-Null
$d.Content.Find.Execute("8+54=62", $true, $false, $false, $false, $false, $true, 1, $false, "86-48=38", 2) | Out-Null
$d.Content.Find.Execute("83-45=38", $true, $false, $false, $false, $false, $true, 1, $false, "19+53=72", 2) | Out-Null
$d.Content.Find.Execute("78-46=32", $true, $false, $false, $false, $false, $true, 1, $false, "88-26=62", 2) | Out-Null
$d.Content.Find.Execute("89+4=93", $true, $false, $false, $false, $false, $true, 1, $false, "70-51=19", 2) | Out-Null
$d.Content.Find.Execute("65-16=49", $true, $false, $false, $false, $false, $true, 1, $false, "15+43=58", 2) | Out-Null
$d.Content.Find.Execute("20+11=31", $true, $false, $false, $false, $false, $true, 1, $false, "45+51=96", 2) | Out-Null
$d.Content.Find.Execute("32-16=16", $true, $false, $false, $false, $false, $true, 1, $false, "52+8=60", 2) | Out-Null
$d.Content.Find.Execute("72-23=49", $true, $false, $false, $false, $false, $true, 1, $false, "70+18=88", 2) | Out-Null
$d.Content.Find.Execute("35-4=31", $true, $false, $false, $false, $false, $true, 1, $false, "82-13=69", 2) | Out-Null
$d.Content.Find.Execute("58-32=26", $true, $false, $false, $false, $false, $true, 1, $false, "37-12=25", 2) | Out-Null
$d.Content.Find.Execute("36+27=63", $true, $false, $false, $false, $false, $true, 1, $false, "58+37=95", 2) | Out-Null
$d.Content.Find.Execute("22+52=74", $true, $false, $false, $false, $false, $true, 1, $false, "41+24=65", 2) | Out-Null
$d.Content.Find.Execute("1+72=73", $true, $false, $false, $false, $false, $true, 1, $false, "88-26=62", 2) | Out-Null
$d.Content.Find.Execute("99-90=9", $true, $false, $false, $false, $false, $true, 1, $false, "7+36=43", 2) | Out-Null
$d.Content.Find.Execute("52-15=37", $true, $false, $false, $false, $false, $true, 1, $false, "97-61=36", 2) | Out-Null
$d.Content.Find.Execute("22+3=25", $true, $false, $false, $false, $false, $true, 1, $false, "69-10=59", 2) | Out-Null
$d.Content.Find.Execute("42-6=36", $true, $false, $false, $false, $false, $true, 1, $false, "64-25=39", 2) | Out-Null
$d.Content.Find.Execute("60+5=65", $true, $false, $false, $false, $false, $true, 1, $false, "90-4=86", 2) | Out-Null
$d.Content.Find.Execute("92+1=93", $true, $false, $false, $false, $false, $true, 1, $false, "29+22=51", 2) | Out-Null
$d.Content.Find.Execute("62+3=65", $true, $false, $false, $false, $false, $true, 1, $false, "86-38=48", 2) | Out-Null
$d.Content.Find.Execute("5+46=51", $true, $false, $false, $false, $false, $true, 1, $false, "64-63=1", 2) | Out-Null
$d.Content.Find.Execute("33-14=19", $true, $false, $false, $false, $false, $true, 1, $false, "32-21=11", 2) | Out-Null
$d.Content.Find.Execute("94-70=24", $true, $false, $false, $false, $false, $true, 1, $false, "66-18=48", 2) | Out-Null
$d.Content.Find.Execute("53+29=82", $true, $false, $false, $false, $false, $true, 1, $false, "90-36=54", 2) | Out-Null
$d.Content.Find.Execute("68-35=33", $true, $false, $false, $false, $false, $true, 1, $false, "95-59=36", 2) | Out-Null
$d.Content.Find.Execute("33+19=52", $true, $false, $false, $false, $false, $true, 1, $false, "88+9=97", 2) | Out-Null
$d.Content.Find.Execute("48-35=13", $true, $false, $false, $false, $false, $true, 1, $false, "43-38=5", 2) | Out-Null
$d.Content.Find.Execute("22+61=83", $true, $false, $false, $false, $false, $true, 1, $false, "29-10=19", 2) | Out-Null
$d.Content.Find.Execute("29-4=25", $true, $false, $false, $false, $false, $true, 1, $false, "72-33=39", 2) | Out-Null
$d.Content.Find.Execute("15+8=23", $true, $false, $false, $false, $false, $true, 1, $false, "48-30=18", 2) | Out-Null
$d.Content.Find.Execute("84-61=23", $true, $false, $false, $false, $false, $true, 1, $false, "28+4=32", 2) | Out-Null
$d.Content.Find.Execute("1+27=28", $true, $false, $false, $false, $false, $true, 1, $false, "57+17=74", 2) | Out-Null
$d.Content.Find.Execute("74-45=29", $true, $false, $false, $false, $false, $true, 1, $false, "16+67=83", 2) | Out-Null
$d.Content.Find.Execute("41+41=82", $true, $false, $false, $false, $false, $true, 1, $false, "55-32=23", 2) | Out-Null
$d.Content.Find.Execute("76-66=10", $true, $false, $false, $false, $false, $true, 1, $false, "60+27=87", 2) | Out-Null
$d.Content.Find.Execute("74+2=76", $true, $false, $false, $false, $false, $true, 1, $false, "50+17=67", 2) | Out-Null
$d.Content.Find.Execute("60-46=14", $true, $false, $false, $false, $false, $true, 1, $false, "73-51=22", 2) | Out-Null
$d.Content.Find.Execute("74+8=82", $true, $false, $false, $false, $false, $true, 1, $false, "29+6=35", 2) | Out-Null
$d.Content.Find.Execute("23+55=78", $true, $false, $false, $false, $false, $true, 1, $false, "72-62=10", 2) | Out-Null
$d.Content.Find.Execute("20-5=15", $true, $false, $false, $false, $false, $true, 1, $false, "2+4=6", 2) | Out-Null
$d.Content.Find.Execute("86+12=98", $true, $false, $false, $false, $false, $true, 1, $false, "20+75=95", 2) | Out-Null
$d.Content.Find.Execute("7+57=64", $true, $false, $false, $false, $false, $true, 1, $false, "79+5=84", 2) | Out-Null
$d.Content.Find.Execute("84-67=17", $true, $false, $false, $false, $false, $true, 1, $false, "68+19=87", 2) | Out-Null
$d.Content.Find.Execute("25+18=43", $true, $false, $false, $false, $false, $true, 1, $false, "13+42=55", 2) | Out-Null
$d.Content.Find.Execute("11+85=96", $true, $false, $false, $false, $false, $true, 1, $false, "10+35=45", 2) | Out-Null
$d.Content.Find.Execute("55-4=51", $true, $false, $false, $false, $false, $true, 1, $false, "41+13=54", 2) | Out-Null
$d.Content.Find.Execute("59+11=70", $true, $false, $false, $false, $false, $true, 1, $false, "73-57=16", 2) | Out-Null
$d.Content.Find.Execute("1+56=57", $true, $false, $false, $false, $false, $true, 1, $false, "33+5=38", 2) | Out-Null
$d.Content.Find.Execute("60-24=36", $true, $false, $false, $false, $false, $true, 1, $false, "69-30=39", 2) | Out-Null
$d.Content.Find.Execute("73-42=31", $true, $false, $false, $false, $false, $true, 1, $false, "95-71=24", 2) | Out-Null
$d.Content.Find.Execute("68-20=48", $true, $false, $false, $false, $false, $true, 1, $false, "37-31=6", 2) | Out-Null
$d.Content.Find.Execute("26+63=89", $true, $false, $false, $false, $false, $true, 1, $false, "72-24=48", 2) | Out-Null
$d.Content.Find.Execute("53-45=8", $true, $false, $false, $false, $false, $true, 1, $false, "99-73=26", 2) | Out-Null
$d.Content.Find.Execute("39+46=85", $true, $false, $false, $false, $false, $true, 1, $false, "51-26=25", 2) | Out-Null
$d.Content.Find.Execute("53+44=97", $true, $false, $false, $false, $false, $true, 1, $false, "27+6=33", 2) | Out-Null
$d.Content.Find.Execute("1+24=25", $true, $false, $false, $false, $false, $true, 1, $false, "26+9=35", 2) | Out-Null
$d.Content.Find.Execute("46-34=12", $true, $false, $false, $false, $false, $true, 1, $false, "58-21=37", 2) | Out-Null
$d.Content.Find.Execute("61+23=84", $true, $false, $false, $false, $false, $true, 1, $false, "44-29=15", 2) | Out-Null
$d.Content.Find.Execute("77-16=61", $true, $false, $false, $false, $false, $true, 1, $false, "26+32=58", 2) | Out-Null
$d.Content.Find.Execute("35+24=59", $true, $false, $false, $false, $false, $true, 1, $false, "0+18=18", 2) | Out-Null
$d.Content.Find.Execute("91-50=41", $true, $false, $false, $false, $false, $true, 1, $false, "29+62=91", 2) | Out-Null
$d.Content.Find.Execute("53+24=77", $true, $false, $false, $false, $false, $true, 1, $false, "14+25=39", 2) | Out-Null
$d.Content.Find.Execute("35+1=36", $true, $false, $false, $false, $false, $true, 1, $false, "7-6=1", 2) | Out-Null
$d.Content.Find.Execute("38+58=96", $true, $false, $false, $false, $false, $true, 1, $false, "94-48=46", 2) | Out-Null
$d.Content.Find.Execute("72+4=76", $true, $false, $false, $false, $false, $true, 1, $false, "25+62=87", 2) | Out-Null
$d.Content.Find.Execute("89-43=46", $true, $false, $false, $false, $false, $true, 1, $false, "9+24=33", 2) | Out-Null
$d.Content.Find.Execute("6+63=69", $true, $false, $false, $false, $false, $true, 1, $false, "72-40=32", 2) | Out-Null
$d.Content.Find.Execute("48-28=20", $true, $false, $false, $false, $false, $true, 1, $false, "67+25=92", 2) | Out-Null
$d.Content.Find.Execute("74-57=17", $true, $false, $false, $false, $false, $true, 1, $false, "24+53=77", 2) | Out-Null
$d.Content.Find.Execute("32+11=43", $true, $false, $false, $false, $false, $true, 1, $false, "72-40=32", 2) | Out-Null
$d.Content.Find.Execute("78-27=51", $true, $false, $false, $false, $false, $true, 1, $false, "39+18=57", 2) | Out-Null
$d.Content.Find.Execute("50-39=11", $true, $false, $false, $false, $false, $true, 1, $false, "89-58=31", 2) | Out-Null
$d.Content.Find.Execute("93-13=80", $true, $false, $false, $false, $false, $true, 1, $false, "28+68=96", 2) | Out-Null
$d.Content.Find.Execute("71-71=0", $true, $false, $false, $false, $false, $true, 1, $false, "77-58=19", 2) | Out-Null
$d.Content.Find.Execute("84-20=64", $true, $false, $false, $false, $false, $true, 1, $false, "90-49=41", 2) | Out-Null
$d.Content.Find.Execute("61+31=92", $true, $false, $false, $false, $false, $true, 1, $false, "24+4=28", 2) | Out-Null
$d.Content.Find.Execute("53+41=94", $true, $false, $false, $false, $false, $true, 1, $false, "45+25=70", 2) | Out-Null
$d.Content.Find.Execute("16+26=42", $true, $false, $false, $false, $false, $true, 1, $false, "37+25=62", 2) | Out-Null
$d.Content.Find.Execute("0+72=72", $true, $false, $false, $false, $false, $true, 1, $false, "55+17=72", 2) | Out-Null
$d.Content.Find.Execute("39+57=96", $true, $false, $false, $false, $false, $true, 1, $false, "52+36=88", 2) | Out-Null
$d.Content.Find.Execute("55-36=19", $true, $false, $false, $false, $false, $true, 1, $false, "98-49=49", 2) | Out-Null
$d.Content.Find.Execute("65-18=47", $true, $false, $false, $false, $false, $true, 1, $false, "32+65=97", 2) | Out-Null
$d.Content.Find.Execute("72-7=65", $true, $false, $false, $false, $false, $true, 1, $false, "80-48=32", 2) | Out-Null
$d.Content.Find.Execute("6+73=79", $true, $false, $false, $false, $false, $true, 1, $false, "11+4=15", 2) | Out-Null
$d.Content.Find.Execute("31+18=49", $true, $false, $false, $false, $false, $true, 1, $false, "60-43=17", 2) | Out-Null
$d.Content.Find.Execute("47+3=50", $true, $false, $false, $false, $false, $true, 1, $false, "44-3=41", 2) | Out-Null
$d.Content.Find.Execute("6+60=66", $true, $false, $false, $false, $false, $true, 1, $false, "27-3=24", 2) | Out-Null
$d.Content.Find.Execute("13+26=39", $true, $false, $false, $false, $false, $true, 1, $false, "67-58=9", 2) | Out-Null
$d.Content.Find.Execute("94-79=15", $true, $false, $false, $false, $false, $true, 1, $false, "40-18=22", 2) | Out-Null
$d.Content.Find.Execute("52-21=31", $true, $false, $false, $false, $false, $true, 1, $false, "27+32=59", 2) | Out-Null
$d.Content.Find.Execute("6+64=70", $true, $false, $false, $false, $false, $true, 1, $false, "5+42=47", 2) | Out-Null
$d.Content.Find.Execute("31+61=92", $true, $false, $false, $false, $false, $true, 1, $false, "53-52=1", 2) | Out-Null
$d.Content.Find.Execute("1+17=18", $true, $false, $false, $false, $false, $true, 1, $false, "95-61=34", 2) | Out-Null
$d.Content.Find.Execute("80-68=12", $true, $false, $false, $false, $false, $true, 1, $false, "82+7=89", 2) | Out-Null
$d.Content.Find.Execute("68-64=4", $true, $false, $false, $false, $false, $true, 1, $false, "9+35=44", 2) | Out-Null
